# Update the cryptos list (Price / Volume(1h) columns) on Sheet1.
# Values that look like plain numbers (e.g. "1.00", "593.20") are written
# with a leading apostrophe so Excel stores them as literal text (matching
# the workbook's existing inline/shared-string cells) instead of coercing
# them into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.182.69"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "2.512.89"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'593.20"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").Value = "'175.35"
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "2.512.87"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "2.939.36"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").Value = "'25.74"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "69.011.39"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "2.511.18"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").Value = "'361.67"
$ws.Range("E19").Value = "  +2.83%  "

$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").Value = "'10.91"
$ws.Range("E21").Value = "  -1.68%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "'70.15"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").Value = "'4.16"
$ws.Range("E25").Value = "  -3.04%  "

$ws.Range("D26").Value = "'8.92"
$ws.Range("E26").Value = "  -3.07%  "

$ws.Range("E27").Value = "  -7.81%  "

$ws.Range("D28").Value = "2.637.76"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "'504.99"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").Value = "0.0₃0880"
$ws.Range("E31").Value = "  -3.60%  "

$ws.Range("D32").Value = "'7.69"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("E34").Value = "  -4.53%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").Value = "'162.64"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("D37").Value = "'0.118"
$ws.Range("E37").Value = "  -4.04%  "

$ws.Range("D38").Value = "'18.69"
$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("D44").Value = "'0.319"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'2.30"
$ws.Range("E45").Value = "  -4.98%  "

$ws.Range("D46").Value = "'149.45"
$ws.Range("E46").Value = "  +2.57%  "

$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").Value = "'0.0736"
$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("E51").Value = "  -1.87%  "
